# "Generate Report for Handoff"
#
# The localization-status report is regenerated: a new handoff UUID
# (1662588e-31fe-447b-83ad-66084ab5ec36) replaces the old one
# (3f104482-1790-487a-b0cf-29f7e2e0ede4) everywhere it is referenced
# (file names / hyperlink display text), the xlf token hashes change
# (18ad89769114cf81ca3e1ee9efd6c1511b0342b5 -> 1c7876ecaa922d2c25869e29f78400515978afa7)
# and the handoff timestamps advance by a few seconds.
#
# Note: the external hyperlink targets (the github blob URLs, i.e. the
# relationship Targets) are left untouched - only the *visible* file
# name / display text is updated, exactly like the source diff shows.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldMd  = "3f104482-1790-487a-b0cf-29f7e2e0ede4.md"
$newMd  = "1662588e-31fe-447b-83ad-66084ab5ec36.md"

$oldZh  = "3f104482-1790-487a-b0cf-29f7e2e0ede4.18ad89769114cf81ca3e1ee9efd6c1511b0342b5.zh-cn.xlf"
$newZh  = "1662588e-31fe-447b-83ad-66084ab5ec36.1c7876ecaa922d2c25869e29f78400515978afa7.zh-cn.xlf"

$oldDe  = "3f104482-1790-487a-b0cf-29f7e2e0ede4.18ad89769114cf81ca3e1ee9efd6c1511b0342b5.de-de.xlf"
$newDe  = "1662588e-31fe-447b-83ad-66084ab5ec36.1c7876ecaa922d2c25869e29f78400515978afa7.de-de.xlf"

$mdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/8c4f6bbe114ce84aa7cc6880da94c944659fcbbc/e2e/$oldMd"
$zhTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/665acd151633c47adee1f9c9056edd8d947decff/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZh"
$deTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/47ddd55d182ac91904c7a0074d015aadebaabf58/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDe"

# ---------------------------------------------------------------------
# Overview sheet: A2 is the handoff markdown file name/link, D2 is the
# latest handoff date.
# ---------------------------------------------------------------------
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdTarget, "", "", $newMd)
$wsOverview.Range("D2").Value = "2016-03-22 02:53:10"

# ---------------------------------------------------------------------
# zh-cn sheet: A2 is the handoff markdown file, D2 is the handoff xlf
# target file, E2 is the handoff datetime.
# ---------------------------------------------------------------------
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdTarget, "", "", $newMd)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $zhTarget, "", "", $newZh)
$wsZhCn.Range("E2").Value = "2016-03-22 02:53:06"

# ---------------------------------------------------------------------
# de-de sheet: A2 is the handoff markdown file, D2 is the handoff xlf
# target file, E2 is the handoff datetime (shares the same timestamp
# string as Overview!D2).
# ---------------------------------------------------------------------
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdTarget, "", "", $newMd)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $deTarget, "", "", $newDe)
$wsDeDe.Range("E2").Value = "2016-03-22 02:53:10"
